$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet: bump Version and Date
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.2.0"
$meta.Range("B8").Value = "2024-03-28T10:46:20+01:00"

# ---------------------------------------------------------------------------
# 2. Mapping Table 0 sheet: append five new concept-map rows
# ---------------------------------------------------------------------------
$map = $wb.Worksheets.Item("Mapping Table 0")

$newRows = @(
    @("49614-1", "Campylobacter sp DNA [Identifier] in Specimen by NAA with probe detection", "related-to", "http://fhir.ch/ig/ch-elm/ValueSet/ch-elm-results-camp-org"),
    @("4992-4", "Campylobacter sp rRNA [Presence] in Specimen by Probe", "related-to", "http://fhir.ch/ig/ch-elm/ValueSet/ch-elm-results-camp-org"),
    @("71429-5", "Campylobacter sp DNA.diarrheagenic [Presence] in Stool by NAA with probe detection", "related-to", "http://fhir.ch/ig/ch-elm/ValueSet/ch-elm-results-camp-diar-org"),
    @("85827-4", "Carbapenem resistance bla OXA-48-like gene [Presence] by Molecular method", "related-to", "http://fhir.ch/ig/ch-elm/ValueSet/ch-elm-results-cpe-org"),
    @("97513-6", "Campylobacter sp [Presence] in Specimen by Organism specific culture", "related-to", "http://fhir.ch/ig/ch-elm/ValueSet/ch-elm-results-camp-org")
)

$startRow = 4
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]

    $map.Cells.Item($r, 1).Value = $values[0]
    $map.Cells.Item($r, 2).Value = $values[1]
    $map.Cells.Item($r, 3).Value = $values[2]
    $map.Cells.Item($r, 4).Value = $values[3]
    # Column E is left blank, matching the other data rows.

    for ($c = 1; $c -le 5; $c++) {
        $cell = $map.Cells.Item($r, $c)
        $cell.WrapText = $true
        $cell.VerticalAlignment = -4160
        $cell.Borders.Item(7).LineStyle = 1
        $cell.Borders.Item(8).LineStyle = 1
        $cell.Borders.Item(9).LineStyle = 1
        $cell.Borders.Item(10).LineStyle = 1
    }
}
